# Update countries & provincias Spain
# - Reorder a few countries (Honduras/Armenia, Bahamas/Islas Feroe/Gibraltar,
#   Islas Malvinas/Groenlandia) together with their refreshed per-country stats.
# - Refresh numeric stats for several other countries (India, Kazajistan,
#   Belgica, Butan) with no reordering.
# - Bump the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 06:34"

# --- Row 6 (India): numeric refresh only ---
$ws.Range("D6").Value = 753050
$ws.Range("E6").Value = 412265

# --- Row 32 (Kazajistan): numeric refresh only ---
$ws.Range("B32").Value = 75153
$ws.Range("C32").Value = 1685
$ws.Range("D32").Value = 46790
$ws.Range("E32").Value = 27778

# --- Row 36 (Belgica): numeric refresh only ---
$ws.Range("B36").Value = 64258
$ws.Range("C36").Value = 164
$ws.Range("D36").Value = 17330
$ws.Range("E36").Value = 37123

# --- Rows 53-54: Honduras now sorts before Armenia, each keeps its own stats ---
$ws.Range("A53").Value = "Honduras"
$ws.Range("B53").Value = 35345
$ws.Range("C53").Value = 734
$ws.Range("D53").Value = 4144
$ws.Range("E53").Value = 30213
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 53
$ws.Range("H53").Value = 988

$ws.Range("A54").Value = "Armenia"
$ws.Range("B54").Value = 35254
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 24206
$ws.Range("E54").Value = 10386
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 662

# --- Rows 176-178: Bahamas now sorts before Islas Feroe / Gibraltar ---
$ws.Range("A176").Value = "Bahamas"
$ws.Range("B176").Value = 194
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 91
$ws.Range("E176").Value = 92
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 11

$ws.Range("A177").Value = "Islas Feroe"
$ws.Range("B177").Value = 191
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 188
$ws.Range("E177").Value = 3
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

$ws.Range("A178").Value = "Gibraltar"
$ws.Range("B178").Value = 180
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 180
$ws.Range("E178").Value = 0
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# --- Row 187 (Butan): numeric refresh only ---
$ws.Range("B187").Value = 92
$ws.Range("C187").Value = 2
$ws.Range("E187").Value = 9

# --- Rows 210-211: Islas Malvinas now sorts before Groenlandia (identical stats) ---
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
